$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 16.64788523577431
    "C2" = 10.96267734130551
    "D2" = 5.974574029107775
    "E2" = 11.4177352427967
    "G2" = 25.50443024336143
    "H2" = 13.43888598846094
    "I2" = 20.307434161262
    "L2" = 9.78230786514214
    "O2" = 20.02624156146513
    "B3" = 15.92148017205632
    "C3" = 10.67135346013609
    "D3" = 5.853827881153969
    "E3" = 11.47562813725585
    "G3" = 25.57791268578199
    "H3" = 13.50029160805024
    "I3" = 20.46081009589008
    "L3" = 9.746698247759586
    "O3" = 20.12239195217228
    "B4" = 15.45818407752698
    "C4" = 10.48780111349629
    "D4" = 5.780192019635623
    "E4" = 11.51341733328288
    "G4" = 25.63569128001132
    "H4" = 13.54096679475966
    "I4" = 20.56049605468548
    "L4" = 9.726398868585127
    "O4" = 20.18763090903874
    "B5" = 15.2652731884131
    "C5" = 10.41191034569298
    "D5" = 5.750358118016663
    "E5" = 11.52938103705639
    "G5" = 25.66239610208187
    "H5" = 13.55828840698893
    "I5" = 20.60250488248532
    "L5" = 9.718526081705591
    "O5" = 20.21576888665874
    "B6" = 15.23299941281097
    "C6" = 10.39924539612471
    "D6" = 5.745416069036995
    "E6" = 11.53206589084682
    "G6" = 25.66702050397793
    "H6" = 13.56120967897997
    "I6" = 20.60956412370606
    "L6" = 9.717243094746337
    "O6" = 20.22053474016746
    "B7" = 15.45559874148448
    "C7" = 10.48678192806613
    "D7" = 5.779788905132284
    "E7" = 11.51363034006839
    "G7" = 25.63603867092751
    "H7" = 13.54119738051457
    "I7" = 20.56105698878224
    "L7" = 9.726291068986413
    "O7" = 20.18800411073858
    "B8" = 16.40113625927705
    "C8" = 10.86324883457127
    "D8" = 5.932864344667042
    "E8" = 11.43723144271516
    "G8" = 25.52712573985021
    "H8" = 13.45944113171967
    "I8" = 20.35917393620955
    "L8" = 9.769708455689912
    "O8" = 20.05810275710667
    "B9" = 18.10961753011954
    "C9" = 11.56096116085301
    "D9" = 6.235013533868945
    "E9" = 11.30519259981688
    "G9" = 25.41494839993782
    "H9" = 13.32274994784225
    "I9" = 20.00703000403804
    "L9" = 9.866994478310261
    "O9" = 19.85288691223534
    "B10" = 19.2664215169903
    "C10" = 12.0447065052338
    "D10" = 6.455584022032326
    "E10" = 11.21900025792681
    "G10" = 25.39545928561434
    "H10" = 13.23680117970775
    "I10" = 19.77499385192553
    "L10" = 9.945483557262905
    "O10" = 19.7327164372205
    "B11" = 19.76976525915967
    "C11" = 12.25774854202919
    "D11" = 6.555098429838517
    "E11" = 11.18213363654282
    "G11" = 25.40043413930478
    "H11" = 13.20086023158403
    "I11" = 19.67523450050843
    "L11" = 9.982622338544367
    "O11" = 19.68477767883743
    "B12" = 19.95697402974259
    "C12" = 12.33735744235986
    "D12" = 6.592622873976101
    "E12" = 11.16850987266825
    "G12" = 25.40431794745306
    "H12" = 13.18770574667322
    "I12" = 19.6382925218134
    "L12" = 9.996883797338745
    "O12" = 19.66759864202051
    "B13" = 19.91680782256674
    "C13" = 12.32026047239952
    "D13" = 6.584549034925284
    "E13" = 11.17142901412007
    "G13" = 25.40339243275173
    "H13" = 13.19051851443403
    "I13" = 19.64621148938187
    "L13" = 9.993803675607788
    "O13" = 19.67125501326094
    "B14" = 19.78523564927496
    "C14" = 12.26431967295601
    "D14" = 6.558189011484283
    "E14" = 11.18100605100814
    "G14" = 25.40071353696176
    "H14" = 13.19976886513524
    "I14" = 19.67217852401479
    "L14" = 9.983791721876758
    "O14" = 19.68334478315681
    "B15" = 19.70419874296824
    "C15" = 12.22991397327263
    "D15" = 6.542020787853358
    "E15" = 11.18691612616618
    "G15" = 25.39933331857534
    "H15" = 13.20549434612146
    "I15" = 19.68819283231893
    "L15" = 9.977684618569652
    "O15" = 19.69087720502629
    "B16" = 19.23305682264563
    "C16" = 12.03063756173654
    "D16" = 6.449060379814471
    "E16" = 11.22145671443709
    "G16" = 25.39541357684919
    "H16" = 13.23921367314543
    "I16" = 19.78163008761403
    "L16" = 9.943084592854287
    "O16" = 19.73598531845599
    "B17" = 18.93808318092274
    "C17" = 11.90654891880058
    "D17" = 6.391791187581547
    "E17" = 11.2432462236474
    "G17" = 25.3965617621867
    "H17" = 13.26070919555718
    "I17" = 19.84043601696034
    "L17" = 9.922219954430989
    "O17" = 19.76538573912954
    "B18" = 18.76627249206931
    "C18" = 11.83451904405347
    "D18" = 6.358776192553818
    "E18" = 11.25599948115133
    "G18" = 25.39852453055402
    "H18" = 13.27336996634341
    "I18" = 19.87480499683342
    "L18" = 9.910354633460354
    "O18" = 19.78292889589853
    "B19" = 18.70773461720719
    "C19" = 11.81001987681348
    "D19" = 6.347586241451893
    "E19" = 11.26035538956398
    "G19" = 25.39941242382312
    "H19" = 13.27770766951149
    "I19" = 19.88653536974773
    "L19" = 9.906360755970514
    "O19" = 19.78897717846606
    "B20" = 18.96970692269758
    "C20" = 11.91982682923605
    "D20" = 6.397895685342299
    "E20" = 11.24090387190212
    "G20" = 25.39630466570533
    "H20" = 13.25839019925947
    "I20" = 19.83411957252758
    "L20" = 9.92442707005433
    "O20" = 19.76219046945
    "B21" = 19.82397451835967
    "C21" = 12.28078018285123
    "D21" = 6.565936233138103
    "E21" = 11.17818390205859
    "G21" = 25.40144605488926
    "H21" = 13.19703943667922
    "I21" = 19.66452871040606
    "L21" = 9.986727174971524
    "O21" = 19.67976722743754
    "B22" = 20.36245287755264
    "C22" = 12.51045109489917
    "D22" = 6.67481341656619
    "E22" = 11.13915604834443
    "G22" = 25.41646558593838
    "H22" = 13.15959929897761
    "I22" = 19.55855744502577
    "L22" = 10.02859279667177
    "O22" = 19.63158107457064
    "B23" = 20.07690203564517
    "C23" = 12.38845901743341
    "D23" = 6.616803194203742
    "E23" = 11.15980631906038
    "G23" = 25.40738020090487
    "H23" = 13.17933826574649
    "I23" = 19.61467054479415
    "L23" = 10.00614603690449
    "O23" = 19.656776776005
    "B24" = 18.95541675571652
    "C24" = 11.91382602974998
    "D24" = 6.395136119650798
    "E24" = 11.24196214459316
    "G24" = 25.39641684208383
    "H24" = 13.25943767465992
    "I24" = 19.83697349067246
    "L24" = 9.923428827401954
    "O24" = 19.76363305636874
    "B25" = 17.66411508794451
    "C25" = 11.37700325996336
    "D25" = 6.153339908382726
    "E25" = 11.33901160308653
    "G25" = 25.43431392481625
    "H25" = 13.3571912005281
    "I25" = 20.09760981660988
    "L25" = 9.839414690061865
    "O25" = 19.90305849767807
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

